# Update the "想去人数" (want-to-go count) column F for several rows on the
# "展览" sheet and the equivalent rows on the "全部类型" sheet.
# Each listed value is incremented by 1.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st sheet) ---
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F4").Value  = 14020
$wsExpo.Range("F5").Value  = 246
$wsExpo.Range("F11").Value = 39
$wsExpo.Range("F17").Value = 15083
$wsExpo.Range("F19").Value = 8450
$wsExpo.Range("F20").Value = 298
$wsExpo.Range("F33").Value = 36
$wsExpo.Range("F40").Value = 239

# --- Sheet "全部类型" (4th sheet) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value  = 14020
$wsAll.Range("F5").Value  = 246
$wsAll.Range("F11").Value = 39
$wsAll.Range("F17").Value = 15083
$wsAll.Range("F19").Value = 8450
$wsAll.Range("F20").Value = 298
$wsAll.Range("F34").Value = 36
$wsAll.Range("F43").Value = 239
